$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the existing header
# formatting used by columns B..H (bold font, centered/top aligned, thin
# box border around each header cell).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

foreach ($addr in "I1", "J1") {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.BorderAround(1, 2)            # xlContinuous, xlThin
}

# New data columns I ("I0") and J ("IF") for rows 2-10; both columns hold
# the same value for any given row.
$values = @{
    2  = 8
    3  = 9
    4  = 8
    5  = 9
    6  = 7
    7  = 7
    8  = 9
    9  = 8
    10 = 7
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 9).Value  = $values[$row]   # column I
    $ws.Cells.Item($row, 10).Value = $values[$row]   # column J
}
